$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D and E columns for data rows so that numeric-looking
# strings (e.g. "1.000", "0.9991") are preserved exactly as text, not coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.154.67"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").Value = "1.823.04"
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "234.86"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("D6").Value = "0.6029"
$ws.Range("E6").Value = "  -3.95%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.07060"
$ws.Range("E8").Value = "  -4.99%  "
$ws.Range("D9").Value = "0.2791"
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("D10").Value = "23.49"
$ws.Range("E10").Value = "  -5.85%  "
$ws.Range("D11").Value = "0.07632"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("D12").Value = "1.823.22"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "4.786"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").Value = "0.6291"
$ws.Range("E14").Value = "  -6.64%  "
$ws.Range("D15").Value = "0.000009938"
$ws.Range("E15").Value = "  -2.52%  "
$ws.Range("D16").Value = "2.064.86"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "78.44"
$ws.Range("E17").Value = "  -3.96%  "
$ws.Range("D18").Value = "5.846"
$ws.Range("E18").Value = "  -5.95%  "
$ws.Range("D19").Value = "29.139.25"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "225.98"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "11.73"
$ws.Range("E22").Value = "  -4.45%  "
$ws.Range("D23").Value = "6.977"
$ws.Range("E23").Value = "  -5.10%  "
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "155.05"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "8.001"
$ws.Range("E26").Value = "  -5.54%  "
$ws.Range("D27").Value = "0.1301"
$ws.Range("E27").Value = "  -3.37%  "
$ws.Range("D28").Value = "16.52"
$ws.Range("E28").Value = "  -4.73%  "
$ws.Range("D29").Value = "1.489"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").Value = "0.06221"
$ws.Range("E30").Value = "  -14.79%  "
$ws.Range("D31").Value = "1.447"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("D32").Value = "3.833"
$ws.Range("E32").Value = "  -4.99%  "
$ws.Range("D33").Value = "3.791"
$ws.Range("E33").Value = "  -6.18%  "
$ws.Range("D34").Value = "1.121"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("D35").Value = "1.737"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("D36").Value = "0.6384"
$ws.Range("E36").Value = "  -8.04%  "
$ws.Range("D37").Value = "2.539"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").Value = "1.213.55"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("E40").Value = "  -5.54%  "
$ws.Range("D41").Value = "6.473"
$ws.Range("E41").Value = "  -6.42%  "
$ws.Range("D42").Value = "0.9067"
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "1.978.63"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "100.47"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("D46").Value = "62.50"
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("D47").Value = "0.00000000116"
$ws.Range("E47").Value = "  -3.05%  "
$ws.Range("D48").Value = "8.507"
$ws.Range("E48").Value = "  -3.94%  "
$ws.Range("D49").Value = "1.593"
$ws.Range("E49").Value = "  -6.31%  "
$ws.Range("D50").Value = "0.4554"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("D51").Value = "0.05502"
$ws.Range("E51").Value = "  -2.82%  "

# Restore original (default) cell style so no extraneous style index is left on the cells,
# matching the source workbook formatting (General format, default style).
$ws.Range("D2:E51").Style = "Normal"
